# daily auto push: 2026-01-20 09:44 UTC
# Insert a new data row at row 689 (shifting all subsequent rows down by one)
# and populate it with the new day's reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(689).Insert()

# Write the date as literal text (matching the sheet's existing convention of
# storing "yyyy/mm/dd" as a string rather than a real date serial) by
# temporarily forcing a text number format, then clearing formatting again so
# the new row's cells carry no explicit style - same as all the other data
# rows in this sheet.
$ws.Cells.Item(689, 1).NumberFormat = "@"
$ws.Cells.Item(689, 1).Value = "2026/01/20"
$ws.Cells.Item(689, 1).ClearFormats()

$ws.Cells.Item(689, 2).Value = "火"
$ws.Cells.Item(689, 3).Value = 16
$ws.Cells.Item(689, 4).Value = 174
